# Generate Report for Handoff
# Refresh the "Latest Handoff Date(time)" values for the most recently
# handed-off file (d0f8de88-3a56-4619-98c5-ea3770e9334b) across the
# Overview summary sheet and each per-locale detail sheet.

$wb = $excel.ActiveWorkbook

# Overview sheet: File Name d0f8de88-... is on row 7 (D = Latest Handoff Date)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D7").Value = "2016-30-20 00:30:49"

# zh-cn detail sheet: same source file is on row 7 (E = Latest Handoff Datetime)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E7").Value = "2016-03-20 00:30:46"

# de-de detail sheet: same source file is on row 7 (E = Latest Handoff Datetime)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E7").Value = "2016-03-20 00:30:49"
